# Update LR-pair data with newly computed TPM values.
# - Rows 2-9 (ECs/FAPs/MuSCs/Resolving-Mac x FAPs/MuSCs) get refreshed numeric values.
# - The old "Neutrophils" cluster rows are replaced by "Resolving-Mac" data (was rows 10-11),
#   and the trailing duplicate rows (old 10-11) are removed, shrinking the table to 9 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old Resolving-Mac rows 10 and 11).
$ws.Rows.Item(11).Delete() | Out-Null
$ws.Rows.Item(10).Delete() | Out-Null

$data = @(
    @("ECs", "Myoc", "Fzd10", "FAPs", 1, 0.5, 0.7898485, 1.579697, 0.02334787276659702, 0.01587825395766005, 2, 0.6666666666666666, 0.039899, 0.119697, 0.4086672402490986, 0.5089958879585649, 0.0315141653015, 0.189084991809, 0.009541510729212293, 0.008081965972410777),
    @("ECs", "Myoc", "Fzd10", "MuSCs", 1, 0.5, 0.7898485, 1.579697, 0.02334787276659702, 0.01587825395766005, 1, 0.5, 0.057733, 0.115466, 0.5913327597509014, 0.4910041120414351, 0.0456003234505, 0.182401293802, 0.01380636203738473, 0.007796287985249277),
    @("FAPs", "Myoc", "Fzd10", "FAPs", 3, 1, 31.76379, 95.29137, 0.9389356661497832, 0.9578169565640426, 2, 0.6666666666666666, 0.039899, 0.119697, 0.4086672402490986, 0.5089958879585649, 1.26734345721, 11.40609111489, 0.3837122474568809, 0.487524892308085),
    @("FAPs", "Myoc", "Fzd10", "MuSCs", 3, 1, 31.76379, 95.29137, 0.9389356661497832, 0.9578169565640426, 1, 0.5, 0.057733, 0.115466, 0.5913327597509014, 0.4910041120414351, 1.83381888807, 11.00291332842, 0.5552234186929023, 0.4702920642559575),
    @("MuSCs", "Myoc", "Fzd10", "FAPs", 2, 1, 1.210782, 2.421564, 0.03579064096986433, 0.02434024256976313, 2, 0.6666666666666666, 0.039899, 0.119697, 0.4086672402490986, 0.5089958879585649, 0.048308991018, 0.289853946108, 0.01462646247190078, 0.01238908337992345),
    @("MuSCs", "Myoc", "Fzd10", "MuSCs", 2, 1, 1.210782, 2.421564, 0.03579064096986433, 0.02434024256976313, 1, 0.5, 0.057733, 0.115466, 0.5913327597509014, 0.4910041120414351, 0.069902077206, 0.279608308824, 0.02116417849796355, 0.01195115918983969),
    @("Resolving-Mac", "Myoc", "Fzd10", "FAPs", 1, 0.3333333333333333, 0.06514966666666668, 0.195449, 0.001925820113755411, 0.001964546908534168, 2, 0.6666666666666666, 0.039899, 0.119697, 0.4086672402490986, 0.5089958879585649, 0.002599406550333333, 0.023394658953, 0.0007870195911046291, 0.0009999462981456024),
    @("Resolving-Mac", "Myoc", "Fzd10", "MuSCs", 1, 0.3333333333333333, 0.06514966666666668, 0.195449, 0.001925820113755411, 0.001964546908534168, 1, 0.5, 0.057733, 0.115466, 0.5913327597509014, 0.4910041120414351, 0.003761285705666667, 0.022567714234, 0.001138800522650782, 0.0009646006103885657)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
